# The site footer block ("Ver no Jupiter Salvar em pdf Salvar em docx",
# the "© 2020 ..." notice, and the blank separator paragraph right before
# them) is dropped from the end of the document, right after the
# "Requisitos" section's last entry (LOQ4088 ...).

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph - the first paragraph of the
# block being removed.
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the 'Ver no Jupiter' paragraph"
}

$targetStart = $findRange.Start
$targetEnd = $findRange.End

# Resolve the paragraph index of the found range within the document's
# Paragraphs collection (Range.Paragraphs clips Range to the found text,
# so we look it up against the full document paragraphs instead).
$count = $d.Paragraphs.Count
$verIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $targetStart -and $p.Range.End -ge $targetEnd) {
        $verIdx = $i
        break
    }
}

if ($verIdx -eq -1) {
    throw "Could not resolve paragraph index for the found range"
}

# The blank separator paragraph right before "Ver no Jupiter ..." and the
# copyright paragraph right after it are both removed together with it.
$startIdx = $verIdx - 1
$endIdx = $verIdx + 1

$startPara = $d.Paragraphs.Item($startIdx)
$endPara = $d.Paragraphs.Item($endIdx)

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()
